# Update cryptocurrency price and volume data per the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.578.05'
$ws.Range("E2").Value = '  -2.25%  '
$ws.Range("D3").Value = '3.148.32'
$ws.Range("E3").Value = '  -4.04%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '527.27'
$ws.Range("E5").Value = '  -4.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.70'
$ws.Range("E6").Value = '  -3.30%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.147.35'
$ws.Range("E8").Value = '  -4.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.444'
$ws.Range("E9").Value = '  -4.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.21'
$ws.Range("E10").Value = '  -7.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.109'
$ws.Range("E11").Value = '  -7.85%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.381'
$ws.Range("E12").Value = '  -5.98%  '
$ws.Range("D13").Value = '3.685.09'
$ws.Range("E13").Value = '  -4.09%  '
$ws.Range("E14").Value = '  -1.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.65'
$ws.Range("E15").Value = '  -4.50%  '
$ws.Range("D16").Value = '3.141.46'
$ws.Range("E16").Value = '  -4.08%  '
$ws.Range("D17").Value = '58.489.86'
$ws.Range("E17").Value = '  -2.56%  '
$ws.Range("E18").Value = '  -6.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.79'
$ws.Range("E19").Value = '  -4.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.06'
$ws.Range("E20").Value = '  -5.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.96'
$ws.Range("E21").Value = '  -6.79%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '345.58'
$ws.Range("E22").Value = '  -7.48%  '
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.509'
$ws.Range("E24").Value = '  -4.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '68.07'
$ws.Range("E25").Value = '  -7.43%  '
$ws.Range("D26").Value = '3.269.60'
$ws.Range("E26").Value = '  -4.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.170'
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("E28").Value = '  -4.63%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.84'
$ws.Range("E30").Value = '  -3.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("E32").Value = '  -7.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.90'
$ws.Range("E33").Value = '  -7.55%  '
$ws.Range("E34").Value = '  -0.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '21.46'
$ws.Range("E35").Value = '  -4.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.83'
$ws.Range("E36").Value = '  -4.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '157.25'
$ws.Range("E37").Value = '  -5.45%  '
$ws.Range("E38").Value = '  -5.55%  '
$ws.Range("E39").Value = '  -9.21%  '
$ws.Range("E40").Value = '  -5.18%  '
$ws.Range("D41").Value = '3.177.73'
$ws.Range("E41").Value = '  -3.98%  '
$ws.Range("E42").Value = '  -5.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.49'
$ws.Range("E43").Value = '  -2.79%  '
$ws.Range("E44").Value = '  -0.66%  '
$ws.Range("E45").Value = '  -6.93%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.92'
$ws.Range("E46").Value = '  -4.44%  '
$ws.Range("E47").Value = '  -0.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.45'
$ws.Range("E48").Value = '  -7.62%  '
$ws.Range("D49").Value = '2.274.29'
$ws.Range("E49").Value = '  -2.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.23'
$ws.Range("E50").Value = '  -2.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.87'
$ws.Range("E51").Value = '  -1.63%  '
